$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column B: B2 mirrors A2 ("TXL") using the plain default style.
$ws.Range("B2").Value = $ws.Range("A2").Value()

# Update/replace the travel note cells in column A.
$ws.Range("A5").Value = "Berlin, TXL"
$ws.Range("A6").Value = "DFA, Düsseldorf"
$ws.Range("A7").Value = "ADF"

# Rows 6 & 7 now carry an explicit (custom) row height.
$ws.Rows.Item(6).RowHeight = 12.8
$ws.Rows.Item(7).RowHeight = 12.8

# Column A grows a bit wider to fit the new text.
$ws.Columns.Item(1).ColumnWidth = 37.6

# Match the author's final selection (bottom pane, cell B2).
$ws.Range("B2").Select() | Out-Null
